$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("email") - shifts loai/trangThai to C/D
$ws.Columns("B:B").Insert()

# Header
$ws.Range("B1").Value = "email"

# Data rows
$ws.Range("B2").Value = "test_email@gmail.com"
$ws.Range("B3").Value = "admin@gmail.com"

# Approximate the original column's width for the new column
$ws.Columns("B:B").ColumnWidth = 19.33

# Hyperlinks for the email addresses (adds the Hyperlink style automatically)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:test_email@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:admin@gmail.com") | Out-Null

# Restore selection to match new cursor position
$ws.Range("E13").Select() | Out-Null
